$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 6.449754000000001
$ws.Range("H2").Value = 19.349262
$ws.Range("I2").Value = 0.03479900749229446
$ws.Range("J2").Value = 0.03479900749229446
$ws.Range("M2").Value = 2.759544333333333
$ws.Range("N2").Value = 8.278632999999999
$ws.Range("O2").Value = 0.2574067337278401
$ws.Range("P2").Value = 0.2574067337278401
$ws.Range("Q2").Value = 17.798382102094
$ws.Range("R2").Value = 160.185438918846
$ws.Range("S2").Value = 0.008957498855562152
$ws.Range("T2").Value = 0.008957498855562152
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 6.449754000000001
$ws.Range("H3").Value = 19.349262
$ws.Range("I3").Value = 0.03479900749229446
$ws.Range("J3").Value = 0.03479900749229446
$ws.Range("O3").Value = 0.6758254232987829
$ws.Range("P3").Value = 0.6758254232987829
$ws.Range("Q3").Value = 46.72993182415801
$ws.Range("R3").Value = 420.569386417422
$ws.Range("S3").Value = 0.02351805396885742
$ws.Range("T3").Value = 0.02351805396885742
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 6.449754000000001
$ws.Range("H4").Value = 19.349262
$ws.Range("I4").Value = 0.03479900749229446
$ws.Range("J4").Value = 0.03479900749229446
$ws.Range("M4").Value = 0.5200313333333334
$ws.Range("N4").Value = 1.560094
$ws.Range("O4").Value = 0.0485078515798926
$ws.Range("P4").Value = 0.0485078515798926
$ws.Range("Q4").Value = 3.354074172292001
$ws.Range("R4").Value = 30.18666755062801
$ws.Range("S4").Value = 0.001688025090563791
$ws.Range("T4").Value = 0.00168802509056379
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 6.449754000000001
$ws.Range("H5").Value = 19.349262
$ws.Range("I5").Value = 0.03479900749229446
$ws.Range("J5").Value = 0.03479900749229446
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1957573333333333
$ws.Range("N5").Value = 0.587272
$ws.Range("O5").Value = 0.01825999139348442
$ws.Range("P5").Value = 0.01825999139348442
$ws.Range("Q5").Value = 1.262586643696
$ws.Range("R5").Value = 11.363279793264
$ws.Range("S5").Value = 0.0006354295773110968
$ws.Range("T5").Value = 0.0006354295773110968
$ws.Range("I6").Value = 0.663783921437469
$ws.Range("J6").Value = 0.6637839214374691
$ws.Range("M6").Value = 2.759544333333333
$ws.Range("N6").Value = 8.278632999999999
$ws.Range("O6").Value = 0.2574067337278401
$ws.Range("P6").Value = 0.2574067337278401
$ws.Range("Q6").Value = 339.5004834429963
$ws.Range("R6").Value = 3055.504350986967
$ws.Range("S6").Value = 0.1708624511182761
$ws.Range("T6").Value = 0.1708624511182761
$ws.Range("I7").Value = 0.663783921437469
$ws.Range("J7").Value = 0.6637839214374691
$ws.Range("O7").Value = 0.6758254232987829
$ws.Range("P7").Value = 0.6758254232987829
$ws.Range("Q7").Value = 891.3638528803909
$ws.Range("R7").Value = 8022.274675923519
$ws.Range("S7").Value = 0.4486020496844035
$ws.Range("T7").Value = 0.4486020496844036
$ws.Range("I8").Value = 0.663783921437469
$ws.Range("J8").Value = 0.6637839214374691
$ws.Range("M8").Value = 0.5200313333333334
$ws.Range("N8").Value = 1.560094
$ws.Range("O8").Value = 0.0485078515798926
$ws.Range("P8").Value = 0.0485078515798926
$ws.Range("Q8").Value = 63.97827602896734
$ws.Range("R8").Value = 575.804484260706
$ws.Range("S8").Value = 0.03219873194220784
$ws.Range("T8").Value = 0.03219873194220784
$ws.Range("I9").Value = 0.663783921437469
$ws.Range("J9").Value = 0.6637839214374691
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.1957573333333333
$ws.Range("N9").Value = 0.587272
$ws.Range("O9").Value = 0.01825999139348442
$ws.Range("P9").Value = 0.01825999139348442
$ws.Range("Q9").Value = 24.08358093812533
$ws.Range("R9").Value = 216.752228443128
$ws.Range("S9").Value = 0.01212068869258152
$ws.Range("T9").Value = 0.01212068869258153
$ws.Range("G10").Value = 55.79038633333334
$ws.Range("H10").Value = 167.371159
$ws.Range("I10").Value = 0.3010114916028843
$ws.Range("J10").Value = 0.3010114916028843
$ws.Range("M10").Value = 2.759544333333333
$ws.Range("N10").Value = 8.278632999999999
$ws.Range("O10").Value = 0.2574067337278401
$ws.Range("P10").Value = 0.2574067337278401
$ws.Range("Q10").Value = 153.9560444606274
$ws.Range("R10").Value = 1385.604400145647
$ws.Range("S10").Value = 0.0774823848680436
$ws.Range("T10").Value = 0.0774823848680436
$ws.Range("G11").Value = 55.79038633333334
$ws.Range("H11").Value = 167.371159
$ws.Range("I11").Value = 0.3010114916028843
$ws.Range("J11").Value = 0.3010114916028843
$ws.Range("O11").Value = 0.6758254232987829
$ws.Range("P11").Value = 0.6758254232987829
$ws.Range("Q11").Value = 404.2140134026977
$ws.Range("R11").Value = 3637.926120624279
$ws.Range("S11").Value = 0.2034312187303173
$ws.Range("T11").Value = 0.2034312187303173
$ws.Range("G12").Value = 55.79038633333334
$ws.Range("H12").Value = 167.371159
$ws.Range("I12").Value = 0.3010114916028843
$ws.Range("J12").Value = 0.3010114916028843
$ws.Range("M12").Value = 0.5200313333333334
$ws.Range("N12").Value = 1.560094
$ws.Range("O12").Value = 0.0485078515798926
$ws.Range("P12").Value = 0.0485078515798926
$ws.Range("Q12").Value = 29.01274899210512
$ws.Range("R12").Value = 261.114740928946
$ws.Range("S12").Value = 0.0146014207585148
$ws.Range("T12").Value = 0.0146014207585148
$ws.Range("G13").Value = 55.79038633333334
$ws.Range("H13").Value = 167.371159
$ws.Range("I13").Value = 0.3010114916028843
$ws.Range("J13").Value = 0.3010114916028843
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.1957573333333333
$ws.Range("N13").Value = 0.587272
$ws.Range("O13").Value = 0.01825999139348442
$ws.Range("P13").Value = 0.01825999139348442
$ws.Range("Q13").Value = 10.92137725424978
$ws.Range("R13").Value = 98.292395288248
$ws.Range("S13").Value = 0.005496467246008575
$ws.Range("T13").Value = 0.005496467246008575
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.07517133333333333
$ws.Range("H14").Value = 0.225514
$ws.Range("I14").Value = 0.0004055794673521549
$ws.Range("J14").Value = 0.000405579467352155
$ws.Range("M14").Value = 2.759544333333333
$ws.Range("N14").Value = 8.278632999999999
$ws.Range("O14").Value = 0.2574067337278401
$ws.Range("P14").Value = 0.2574067337278401
$ws.Range("Q14").Value = 0.2074386269291111
$ws.Range("R14").Value = 1.866947642362
$ws.Range("S14").Value = 0.0001043988859581953
$ws.Range("T14").Value = 0.0001043988859581954
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.07517133333333333
$ws.Range("H15").Value = 0.225514
$ws.Range("I15").Value = 0.0004055794673521549
$ws.Range("J15").Value = 0.000405579467352155
$ws.Range("O15").Value = 0.6758254232987829
$ws.Range("P15").Value = 0.6758254232987829
$ws.Range("Q15").Value = 0.5446333738926666
$ws.Range("R15").Value = 4.901700365033999
$ws.Range("S15").Value = 0.000274100915204565
$ws.Range("T15").Value = 0.000274100915204565
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.07517133333333333
$ws.Range("H16").Value = 0.225514
$ws.Range("I16").Value = 0.0004055794673521549
$ws.Range("J16").Value = 0.000405579467352155
$ws.Range("M16").Value = 0.5200313333333334
$ws.Range("N16").Value = 1.560094
$ws.Range("O16").Value = 0.0485078515798926
$ws.Range("P16").Value = 0.0485078515798926
$ws.Range("Q16").Value = 0.03909144870177778
$ws.Range("R16").Value = 0.351823038316
$ws.Range("S16").Value = 0.00001967378860617023
$ws.Range("T16").Value = 0.00001967378860617023
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.07517133333333333
$ws.Range("H17").Value = 0.225514
$ws.Range("I17").Value = 0.0004055794673521549
$ws.Range("J17").Value = 0.000405579467352155
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.1957573333333333
$ws.Range("N17").Value = 0.587272
$ws.Range("O17").Value = 0.01825999139348442
$ws.Range("P17").Value = 0.01825999139348442
$ws.Range("Q17").Value = 0.01471533975644444
$ws.Range("R17").Value = 0.132438057808
$ws.Range("S17").Value = 0.000007405877583224345
$ws.Range("T17").Value = 0.000007405877583224345
